$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 472
$ws1.Range("F3").Value = 5671
$ws1.Range("F6").Value = 92
$ws1.Range("F8").Value = 53
$ws1.Range("F9").Value = 536

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 472
$ws4.Range("F3").Value = 5671
$ws4.Range("F7").Value = 92
$ws4.Range("F10").Value = 53
$ws4.Range("F11").Value = 536
